$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 4 with a new skill entry (FieryFlicker), matching the style
# of the existing data rows (copy format from row 3).
$ws.Range("A3:E3").Copy() | Out-Null
$ws.Range("A4:E4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "FieryFlicker"
$ws.Cells.Item(4, 3).Value = "DamageSkill"
$ws.Cells.Item(4, 4).Value = 5
$ws.Cells.Item(4, 5).Value = 2

# Update the active selection as recorded in the sheet view.
$ws.Range("E9").Select() | Out-Null
